# Apply scraped crypto price/volume updates (and one B/C/D/E row swap)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '42.716.61'
$ws.Range('E2').Value = '  +0.50%  '

# Row 3
$ws.Range('D3').Value = '2.305.48'
$ws.Range('E3').Value = '  +0.06%  '

# Row 4
$ws.Range('E4').Value = '  -0.12%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '309.98'
$ws.Range('E5').Value = '  -2.09%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '105.95'
$ws.Range('E6').Value = '  +1.98%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.621'
$ws.Range('E7').Value = '  -1.38%  '

# Row 8
$ws.Range('E8').Value = '  -0.07%  '

# Row 9
$ws.Range('E9').Value = '  -0.12%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.72'
$ws.Range('E10').Value = '  -0.81%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0908'
$ws.Range('E11').Value = '  +0.28%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.29'
$ws.Range('E12').Value = '  -2.96%  '

# Row 13
$ws.Range('E13').Value = '  -0.18%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.992'
$ws.Range('E14').Value = '  -0.62%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.30'
$ws.Range('E15').Value = '  -0.21%  '

# Row 16
$ws.Range('D16').Value = '2.656.14'
$ws.Range('E16').Value = '  +0.09%  '

# Row 17
$ws.Range('D17').Value = '2.304.03'
$ws.Range('E17').Value = '  +0.25%  '

# Row 18
$ws.Range('D18').Value = '42.692.48'
$ws.Range('E18').Value = '  +0.24%  '

# Row 19
$ws.Range('E19').Value = '  -3.79%  '

# Row 20
$ws.Range('E20').Value = '  -0.93%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.37'
$ws.Range('E21').Value = '  -2.11%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '73.63'
$ws.Range('E22').Value = '  -0.49%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.44'
$ws.Range('E23').Value = '  -2.80%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '267.13'
$ws.Range('E24').Value = '  -0.05%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.24'
$ws.Range('E25').Value = '  -0.02%  '

# Row 26
$ws.Range('E26').Value = '  +0.20%  '

# Row 27
$ws.Range('B27').Value = 'Filecoin'
$ws.Range('C27').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.36'
$ws.Range('E27').Value = '  +11.38%  '

# Row 28
$ws.Range('B28').Value = 'Cosmos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.96'
$ws.Range('E28').Value = '  +0.36%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.30'
$ws.Range('E29').Value = '  +1.48%  '

# Row 30
$ws.Range('E30').Value = '  -1.31%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '37.75'
$ws.Range('E31').Value = '  -1.05%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '165.02'
$ws.Range('E32').Value = '  -0.30%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0864'
$ws.Range('E33').Value = '  -2.14%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.82'
$ws.Range('E34').Value = '  +6.35%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.130'
$ws.Range('E35').Value = '  -1.07%  '

# Row 36
$ws.Range('E36').Value = '  -0.97%  '

# Row 37
$ws.Range('E37').Value = '  +0.47%  '

# Row 38
$ws.Range('E38').Value = '  -1.02%  '

# Row 39
$ws.Range('E39').Value = '  +2.32%  '

# Row 40
$ws.Range('E40').Value = '  -1.72%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '108.07'
$ws.Range('E41').Value = '  +10.26%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.58'
$ws.Range('E42').Value = '  -3.77%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '72.04'
$ws.Range('E43').Value = '  +2.70%  '

# Row 44
$ws.Range('E44').Value = '  +1.24%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.00'
$ws.Range('E45').Value = '  -0.28%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '12.37'
$ws.Range('E46').Value = '  +0.19%  '

# Row 47
$ws.Range('D47').Value = '1.721.50'
$ws.Range('E47').Value = '  +5.04%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '111.40'
$ws.Range('E48').Value = '  -4.89%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '76.58'
$ws.Range('E49').Value = '  -4.78%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.90'
$ws.Range('E50').Value = '  +0.29%  '

# Row 51
$ws.Range('E51').Value = '  -2.33%  '
